$d = $word.ActiveDocument

# --- 1) Add justify (jc=both) to four existing paragraphs identified by their text ---
function Set-JustifyByText($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $par = $doc.Paragraphs.Item($i)
        if ($par.Range.Text -like "*$needle*") {
            $par.Alignment = 3  # wdAlignParagraphJustify
            return $i
        }
    }
    return -1
}

Set-JustifyByText $d "respecto a criptomonedas" | Out-Null
Set-JustifyByText $d "debe actualizarse de" | Out-Null
Set-JustifyByText $d "depositar MXN" | Out-Null
Set-JustifyByText $d "activando el" | Out-Null

# --- 2) Insert the new "Sprint 7" block right after the "Crear el metodo para depositar MXN..." paragraph ---
$depositParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*depositar MXN*") {
        $depositParaIndex = $i
        break
    }
}

$depositPara = $d.Paragraphs.Item($depositParaIndex)
$depositPara.Range.InsertParagraphAfter() | Out-Null
$newPara = $d.Paragraphs.Item($depositParaIndex + 1)

$sprint7Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Sprint </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>7</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>–</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Gestionar </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>transacciones de depósitos</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="360"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> Semanas) </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>05</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>febrero</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> – </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>11</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> febrero 2025</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">Cada vez que se registra un depósito debe ser incluido en el historial de transacciones, incluyendo el id de la transacción, id del usuario, id de la criptomoneda, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>el tipo (deposito, retiro, compra o venta),</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> fecha de la transacción, monto en MXN, cantidad de la criptomoneda. </w:t></w:r></w:p>
'@
$newPara.Range.InsertXML($sprint7Xml) | Out-Null

# --- 3) Replace the final (empty) paragraph with the new closing content ---
$lastPara = $d.Paragraphs.Last
$tailXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Crear una entidad con id compuesto </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Se creó una entidad que maneja un id compuesto por los </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ids</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de otras dos entidades. Esto se realizó para asegurar que dos objetos del mismo tipo no pueden contener las mismas </w:t></w:r><w:r><w:t xml:space="preserve">claves de las </w:t></w:r><w:r><w:t xml:space="preserve">entidades juntas. Es decir, solo un objeto puede tener esta combinación de claves identificadoras. </w:t></w:r><w:r><w:t>De esta manera, las billeteras de criptomonedas tendrán una</w:t></w:r><w:r><w:t xml:space="preserve"> clave compuesta,</w:t></w:r><w:r><w:t xml:space="preserve"> combina</w:t></w:r><w:r><w:t>ndo</w:t></w:r><w:r><w:t xml:space="preserve"> el id del usuario y el id de la criptomoneda, evitando que el mismo usuario pueda tener dos billeteras de la misma criptomoneda. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p>
'@
$lastPara.Range.InsertXML($tailXml) | Out-Null

Write-Output "Done. Paragraphs.Count = $($d.Paragraphs.Count)"
